$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = '83-38='
$t.Cell(1, 2).Range.Text = '3+69='
$t.Cell(1, 3).Range.Text = '58+19='
$t.Cell(1, 4).Range.Text = '18+75='
$t.Cell(1, 5).Range.Text = '42-35='
$t.Cell(2, 1).Range.Text = '69+28='
$t.Cell(2, 2).Range.Text = '7+37='
$t.Cell(2, 3).Range.Text = '84-37='
$t.Cell(2, 4).Range.Text = '83-9='
$t.Cell(2, 5).Range.Text = '81-6='
$t.Cell(3, 1).Range.Text = '72-4='
$t.Cell(3, 2).Range.Text = '18+27='
$t.Cell(3, 3).Range.Text = '44+29='
$t.Cell(3, 4).Range.Text = '68+9='
$t.Cell(3, 5).Range.Text = '9+55='
$t.Cell(4, 1).Range.Text = '28+47='
$t.Cell(4, 2).Range.Text = '35-26='
$t.Cell(4, 3).Range.Text = '90-45='
$t.Cell(4, 4).Range.Text = '80-8='
$t.Cell(4, 5).Range.Text = '18+59='
$t.Cell(5, 1).Range.Text = '8+86='
$t.Cell(5, 2).Range.Text = '23-14='
$t.Cell(5, 3).Range.Text = '38+5='
$t.Cell(5, 4).Range.Text = '42-19='
$t.Cell(5, 5).Range.Text = '93-27='
$t.Cell(6, 1).Range.Text = '22-16='
$t.Cell(6, 2).Range.Text = '93-28='
$t.Cell(6, 3).Range.Text = '81-26='
$t.Cell(6, 4).Range.Text = '57+18='
$t.Cell(6, 5).Range.Text = '45+18='
$t.Cell(7, 1).Range.Text = '29+46='
$t.Cell(7, 2).Range.Text = '9+77='
$t.Cell(7, 3).Range.Text = '62-17='
$t.Cell(7, 4).Range.Text = '34-18='
$t.Cell(7, 5).Range.Text = '44-28='
$t.Cell(8, 1).Range.Text = '92-7='
$t.Cell(8, 2).Range.Text = '90-88='
$t.Cell(8, 3).Range.Text = '76+8='
$t.Cell(8, 4).Range.Text = '18+13='
$t.Cell(8, 5).Range.Text = '77-39='
$t.Cell(9, 1).Range.Text = '22+39='
$t.Cell(9, 2).Range.Text = '17+75='
$t.Cell(9, 3).Range.Text = '10-9='
$t.Cell(9, 4).Range.Text = '73-39='
$t.Cell(9, 5).Range.Text = '55-26='
$t.Cell(10, 1).Range.Text = '34-27='
$t.Cell(10, 2).Range.Text = '58+18='
$t.Cell(10, 3).Range.Text = '80-35='
$t.Cell(10, 4).Range.Text = '3+29='
$t.Cell(10, 5).Range.Text = '4+48='
$t.Cell(11, 1).Range.Text = '53+29='
$t.Cell(11, 2).Range.Text = '53-28='
$t.Cell(11, 3).Range.Text = '7+17='
$t.Cell(11, 4).Range.Text = '83-17='
$t.Cell(11, 5).Range.Text = '54-27='
$t.Cell(12, 1).Range.Text = '48+29='
$t.Cell(12, 2).Range.Text = '91-43='
$t.Cell(12, 3).Range.Text = '34+49='
$t.Cell(12, 4).Range.Text = '30-11='
$t.Cell(12, 5).Range.Text = '49+29='
$t.Cell(13, 1).Range.Text = '82-78='
$t.Cell(13, 2).Range.Text = '67-48='
$t.Cell(13, 3).Range.Text = '18+39='
$t.Cell(13, 4).Range.Text = '46+28='
$t.Cell(13, 5).Range.Text = '58+5='
$t.Cell(14, 1).Range.Text = '50-41='
$t.Cell(14, 2).Range.Text = '70-51='
$t.Cell(14, 3).Range.Text = '65+18='
$t.Cell(14, 4).Range.Text = '85-77='
$t.Cell(14, 5).Range.Text = '4+19='
$t.Cell(15, 1).Range.Text = '28+4='
$t.Cell(15, 2).Range.Text = '84-56='
$t.Cell(15, 3).Range.Text = '87-79='
$t.Cell(15, 4).Range.Text = '6+36='
$t.Cell(15, 5).Range.Text = '84-7='
$t.Cell(16, 1).Range.Text = '81-17='
$t.Cell(16, 2).Range.Text = '8+9='
$t.Cell(16, 3).Range.Text = '27+47='
$t.Cell(16, 4).Range.Text = '93-67='
$t.Cell(16, 5).Range.Text = '67-19='
$t.Cell(17, 1).Range.Text = '9+53='
$t.Cell(17, 2).Range.Text = '93-34='
$t.Cell(17, 3).Range.Text = '29+39='
$t.Cell(17, 4).Range.Text = '35-28='
$t.Cell(17, 5).Range.Text = '90-57='
$t.Cell(18, 1).Range.Text = '29+6='
$t.Cell(18, 2).Range.Text = '62-24='
$t.Cell(18, 3).Range.Text = '83-55='
$t.Cell(18, 4).Range.Text = '67+25='
$t.Cell(18, 5).Range.Text = '41-3='
$t.Cell(19, 1).Range.Text = '91-48='
$t.Cell(19, 2).Range.Text = '23+28='
$t.Cell(19, 3).Range.Text = '15+27='
$t.Cell(19, 4).Range.Text = '12+49='
$t.Cell(19, 5).Range.Text = '57+6='
$t.Cell(20, 1).Range.Text = '39+47='
$t.Cell(20, 2).Range.Text = '22+59='
$t.Cell(20, 3).Range.Text = '63-48='
$t.Cell(20, 4).Range.Text = '77-59='
$t.Cell(20, 5).Range.Text = '93-27='
